# g1.9 - ajuste na coluna período para não inserir valor dinâmico
#
# The "Período" column (C) previously held literal year ranges
# ("2012 / 2022" and "2021 / 2022"). These are replaced with relative,
# non-dynamic labels ("atual/dez anos antes" and "atual/ano anterior").
# Along with that, the two blocks of rows (the "ten years ago" block and
# the "year before" block) swap position, so that the "ano anterior"
# (year before) comparison is listed first (rows 2-7) followed by the
# "dez anos antes" (ten years before) comparison (rows 8-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-13: Atividade, Variação, Período
$data = @(
    @("Transporte, armazenagem e correio", 17.77773366884247, "atual/ano anterior"),
    @("Agropecuária", 6.521894003851491, "atual/ano anterior"),
    @("Construção", 6.438360842135382, "atual/ano anterior"),
    @("Informação e comunicação", 4.037854550070413, "atual/ano anterior"),
    @("Atividades imobiliárias", 3.551068180916705, "atual/ano anterior"),
    @("Administração, defesa, educação e saúde públicas e seguridade social", 1.274864519824148, "atual/ano anterior"),
    @("Atividades financeiras, de seguros e serviços relacionados", 35.19143837025082, "atual/dez anos antes"),
    @("Atividades imobiliárias", 31.76377306549257, "atual/dez anos antes"),
    @("Informação e comunicação", 21.77658698762782, "atual/dez anos antes"),
    @("Agropecuária", 9.730306427073359, "atual/dez anos antes"),
    @("Administração, defesa, educação e saúde públicas e seguridade social", 4.760483082368495, "atual/dez anos antes"),
    @("Eletricidade e gás, água, esgoto, atividades de gestão de resíduos e descontaminação", 2.412539862254022, "atual/dez anos antes")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
